$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force the data range to Text format so that numeric-looking
# strings (e.g. "24.445.10", "1.000") are preserved verbatim as text,
# matching the original inline-string cell contents exactly.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "24.445.10"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").Value = "1.663.04"
$ws.Range("E3").Value = "  -2.33%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "307.53"
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "0.3612"
$ws.Range("E7").Value = "  -3.03%  "
$ws.Range("D8").Value = "47.57"
$ws.Range("E8").Value = "  -1.99%  "
$ws.Range("D9").Value = "0.3257"
$ws.Range("E9").Value = "  -5.26%  "
$ws.Range("D10").Value = "1.127"
$ws.Range("E10").Value = "  -5.97%  "
$ws.Range("D11").Value = "0.07031"
$ws.Range("E11").Value = "  -5.72%  "
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").Value = "5.930"
$ws.Range("E13").Value = "  -4.90%  "
$ws.Range("D14").Value = "19.45"
$ws.Range("E14").Value = "  -6.77%  "
$ws.Range("D15").Value = "1.666.09"
$ws.Range("E15").Value = "  -2.38%  "
$ws.Range("D16").Value = "6.588"
$ws.Range("E16").Value = "  -5.52%  "
$ws.Range("D17").Value = "0.00001050"
$ws.Range("E17").Value = "  -6.56%  "
$ws.Range("D18").Value = "0.06573"
$ws.Range("E18").Value = "  -2.21%  "
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").Value = "77.32"
$ws.Range("E20").Value = "  -7.70%  "
$ws.Range("D21").Value = "5.951"
$ws.Range("E21").Value = "  -5.97%  "
$ws.Range("D22").Value = "15.73"
$ws.Range("E22").Value = "  -8.41%  "
$ws.Range("D23").Value = "12.49"
$ws.Range("E23").Value = "  -3.98%  "
$ws.Range("D24").Value = "24.438.89"
$ws.Range("E24").Value = "  -0.85%  "
$ws.Range("D25").Value = "2.458"
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("D26").Value = "2.338"
$ws.Range("E26").Value = "  -15.05%  "
$ws.Range("D27").Value = "147.25"
$ws.Range("E27").Value = "  -2.26%  "
$ws.Range("D28").Value = "18.55"
$ws.Range("E28").Value = "  -8.24%  "
$ws.Range("D29").Value = "1.848.14"
$ws.Range("E29").Value = "  -2.44%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "1.206"
$ws.Range("E30").Value = "  +2.53%  "
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").Value = "124.56"
$ws.Range("E31").Value = "  -5.01%  "
$ws.Range("D32").Value = "4.014"
$ws.Range("E32").Value = "  -4.36%  "
$ws.Range("D33").Value = "5.718"
$ws.Range("E33").Value = "  -15.04%  "
$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").Value = "0.08408"
$ws.Range("E34").Value = "  -4.47%  "
$ws.Range("B35").Value = "WEMIXTOKEN"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "1.694"
$ws.Range("E35").Value = "  -5.29%  "
$ws.Range("D36").Value = "12.39"
$ws.Range("E36").Value = "  -8.49%  "
$ws.Range("D37").Value = "5.215"
$ws.Range("E37").Value = "  -5.50%  "
$ws.Range("D38").Value = "0.06058"
$ws.Range("E38").Value = "  -8.30%  "
$ws.Range("D39").Value = "0.02217"
$ws.Range("E39").Value = "  -7.82%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "8.263"
$ws.Range("E40").Value = "  -8.42%  "
$ws.Range("D41").Value = "1.208"
$ws.Range("E41").Value = "  -5.00%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "0.2059"
$ws.Range("E42").Value = "  -7.48%  "
$ws.Range("D43").Value = "0.9998"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").Value = "0.5926"
$ws.Range("E44").Value = "  -7.61%  "
$ws.Range("D45").Value = "3.744"
$ws.Range("E45").Value = "  -1.87%  "
$ws.Range("D46").Value = "12.71"
$ws.Range("E46").Value = "  -8.54%  "
$ws.Range("D47").Value = "0.5639"
$ws.Range("E47").Value = "  -7.82%  "
$ws.Range("D48").Value = "122.16"
$ws.Range("E48").Value = "  -5.80%  "
$ws.Range("D49").Value = "1.942"
$ws.Range("E49").Value = "  -8.14%  "
$ws.Range("D50").Value = "0.06915"
$ws.Range("E50").Value = "  -4.49%  "
$ws.Range("D51").Value = "74.53"
$ws.Range("E51").Value = "  -6.04%  "

# Restore the default (Normal) style so no stray style index lingers on
# the cells themselves (keeps cell formatting identical to the original).
$dataRange.Style = "Normal"
